# Add a "Percentage" column (C) next to the existing "Certificate" (A) and
# "Frequency" (B) columns, computing each row's share of the total frequency.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column, matching the style used by the existing headers (A1/B1).
# Copy direct formatting (bold font, border, alignment) from B1 since it is
# applied as a cell format rather than a named style.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C1").Value = "Percentage"

# Determine how many data rows exist by looking at column B (Frequency).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

# Compute the total of the Frequency column.
$total = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $total += $ws.Cells.Item($r, 2).Value2
}

# Fill each row's percentage as text, formatted like "25.67%".
# Force a text number format first so Excel stores the literal string
# instead of re-interpreting it as a percentage number, then restore the
# default (Normal) style so the cell ends up with no special formatting,
# matching the plain data cells in column B.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $freq = $ws.Cells.Item($r, 2).Value2
    $pct = [Math]::Round(($freq / $total) * 100, 2)
    $text = "{0:N2}%" -f $pct
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $ws.Cells.Item($r, 2).Style
}

$wb.Save()
